$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update invoice header text values (A1:A5) - fixed incorrect dates / generic demo data
$ws.Range("A1").Value = "Ivanov Ivan Invoice"
$ws.Range("A2").Value = "Invoice date: October 1, 2000"
$ws.Range("A3").Value = "Contract: dated as of September 1, 2000"
$ws.Range("A4").Value = "Invoice number: 2000-10-II"
$ws.Range("A5").Value = "Date of service: October 2000"

# Update bank address and contact tel # (SWIFT code removed / address genericized)
$ws.Range("B21").Value = "1 Lenina str., Moscow, 1000000, tel +7 495 755-58-58, SWIFT "

# Update "your address" value cell - zip code changed from 650000 to 1000000
$prefix = "PR. LENINA, D. 1, KV. 1, MOSCOW, RUSSIA, "
$suffix = "1000000"
$ws.Range("B29").Value = "$prefix$suffix"

# Give the new zip-code portion its own text run/font so it is stored as
# rich text (mirrors how the source workbook represents this cell)
$suffixChars = $ws.Range("B29").Characters($prefix.Length + 1, $suffix.Length)
$suffixChars.Font.Name = "Arial"
$suffixChars.Font.Size = 10
$suffixChars.Font.ColorIndex = -4105

# Update view selection state (cursor moved from B19 to A21)
$ws.Range("A21").Select()
